# Add a new "2022" column (S) to the transport statistics table, mirroring
# the formatting of the existing "2021" column (R) / "2020" column (Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the per-row formatting from column R (identical to column Q) into
# column S for every data row so the new column visually matches its
# neighbours (number format, font, borders, etc.).
$ws.Range("R3:R15").Copy($ws.Range("S3:S15"))

# Header row: 2022
$ws.Range("S3").Value = 2022

# Data rows (values taken from the updated source table)
$ws.Range("S4").Value = 10444.200000000001
$ws.Range("S5").Value = 21.7
$ws.Range("S6").Value = 7361.6
$ws.Range("S7").Value = 143.1
$ws.Range("S8").Value = 844.2

# Row 9 has no reported value for 2022 yet - leave it blank (format only).
$ws.Range("S9").Value = ""

# Rows 10-12 are recorded as formatted text in the source table.
$ws.Range("S10").Value = "2 756,0"
$ws.Range("S11").Value = "1 013,8"
$ws.Range("S12").Value = "1 451,1"

$ws.Range("S13").Value = 273.39999999999998

# Row 14 has no data ("-") for 2022.
$ws.Range("S14").Value = "-"

$ws.Range("S15").Value = 17.7

# Match the author's final selection from the source commit.
$ws.Range("T3").Select()
